# Insert a new weekly price record for "Comercializadora del Agro de Limarí - Haba"
# A new row is inserted at row 102 (pushing existing rows 102-118 down to 103-119),
# and populated with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("102:102").Insert()

$ws.Range("A102").Value = 2
$ws.Range("B102").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 45218
$ws.Range("E102").Value = 4
$ws.Range("F102").Value = 100112026
$ws.Range("G102").Value = "Haba"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 8000
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = 9000
$ws.Range("N102").Value = "`$/saco 25 kilos"
$ws.Range("O102").Value = "Provincia de Limarí"
$ws.Range("P102").Value = 360
$ws.Range("Q102").Value = 25
$ws.Range("R102").Value = "Hortaliza"
